$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MANU")

$ws.Range("D8").Value = 769600
$ws.Range("E8").Value = 758100
$ws.Range("F8").Value = 672200
$ws.Range("G8").Value = 515400
$ws.Range("H8").Value = 565000
$ws.Range("I8").Value = 473700
$ws.Range("J8").Value = 417800
$ws.Range("D14").Value = 2500
$ws.Range("E14").Value = -6200
$ws.Range("F14").Value = 19700
$ws.Range("G14").Value = 3000
$ws.Range("I14").Value = 8100
$ws.Range("J14").Value = 14000
$ws.Range("D15").Value = 194500
$ws.Range("E15").Value = 175600
$ws.Range("F15").Value = 127900
$ws.Range("G15").Value = 143500
$ws.Range("H15").Value = 83400
$ws.Range("I15").Value = 64500
$ws.Range("J15").Value = 59700
$ws.Range("D17").Value = 712000
$ws.Range("E17").Value = 652700
$ws.Range("F17").Value = 582400
$ws.Range("G17").Value = 474200
$ws.Range("H17").Value = 476400
$ws.Range("I17").Value = 392800
$ws.Range("J17").Value = 359300
$ws.Range("D18").Value = 57600
$ws.Range("E18").Value = 105400
$ws.Range("F18").Value = 89800
$ws.Range("G18").Value = 41300
$ws.Range("H18").Value = 88600
$ws.Range("I18").Value = 80900
$ws.Range("J18").Value = 58500
$ws.Range("D20").Value = 2100
$ws.Range("E20").Value = -4400
$ws.Range("G20").Value = -8300
$ws.Range("I20").Value = -30800
$ws.Range("J20").Value = -9400
$ws.Range("D21").Value = 256400
$ws.Range("E21").Value = 278600
$ws.Range("F21").Value = 218100
$ws.Range("G21").Value = 178100
$ws.Range("H21").Value = 168100
$ws.Range("I21").Value = 115300
$ws.Range("J21").Value = 109400
$ws.Range("D22").Value = 25600
$ws.Range("E22").Value = 27300
$ws.Range("F22").Value = 25100
$ws.Range("G22").Value = 37700
$ws.Range("H22").Value = 30900
$ws.Range("I22").Value = 61500
$ws.Range("J22").Value = 55200
$ws.Range("D23").Value = 34000
$ws.Range("E23").Value = 73700
$ws.Range("F23").Value = 63700
$ws.Range("H23").Value = 52800
$ws.Range("I23").Value = -11500
$ws.Range("D24").Value = 19000
$ws.Range("E24").Value = 22600
$ws.Range("F24").Value = 16300
$ws.Range("H24").Value = 21700
$ws.Range("I24").Value = -202400
$ws.Range("J24").Value = -36500
$ws.Range("D26").Value = 15100
$ws.Range("E26").Value = 51100
$ws.Range("F26").Value = 47400
$ws.Range("H26").Value = 31100
$ws.Range("I26").Value = 191000
$ws.Range("J26").Value = 30400
$ws.Range("D27").Value = 15100
$ws.Range("E27").Value = 51100
$ws.Range("F27").Value = 47400
$ws.Range("H27").Value = 31100
$ws.Range("I27").Value = 190800
$ws.Range("J27").Value = 30000
$ws.Range("D29").Value = -63700
$ws.Range("D32").Value = -2100
$ws.Range("E32").Value = 4400
$ws.Range("G32").Value = 8300
$ws.Range("I32").Value = 30800
$ws.Range("J32").Value = 9400
$ws.Range("D33").Value = -48600
$ws.Range("E33").Value = 51100
$ws.Range("F33").Value = 47400
$ws.Range("H33").Value = 31100
$ws.Range("I33").Value = 190800
$ws.Range("J33").Value = 30000
$ws.Range("D35").Value = -48600
$ws.Range("E35").Value = 51100
$ws.Range("F35").Value = 47400
$ws.Range("H35").Value = 31100
$ws.Range("I35").Value = 190800
$ws.Range("J35").Value = 30000
$ws.Range("D41").Value = 315700
$ws.Range("E41").Value = 378600
$ws.Range("F41").Value = 298900
$ws.Range("G41").Value = 203100
$ws.Range("H41").Value = 86600
$ws.Range("I41").Value = 123200
$ws.Range("J41").Value = 184200
$ws.Range("D43").Value = 220200
$ws.Range("E43").Value = 135300
$ws.Range("F43").Value = 167800
$ws.Range("G43").Value = 109200
$ws.Range("H43").Value = 163200
$ws.Range("I43").Value = 89500
$ws.Range("J43").Value = 100000
$ws.Range("D44").Value = 1800
$ws.Range("E44").Value = 2100
$ws.Range("F45").Value = 10300
$ws.Range("D46").Value = 539300
$ws.Range("E46").Value = 520200
$ws.Range("F46").Value = 478200
$ws.Range("G46").Value = 312400
$ws.Range("H46").Value = 249800
$ws.Range("I46").Value = 213000
$ws.Range("J46").Value = 193300
$ws.Range("E47").Value = 20100
$ws.Range("F47").Value = 14600
$ws.Range("G47").Value = 5000
$ws.Range("J47").Value = 3900
$ws.Range("D48").Value = 356200
$ws.Range("E48").Value = 337400
$ws.Range("F48").Value = 338000
$ws.Range("G48").Value = 344600
$ws.Range("H48").Value = 350200
$ws.Range("I48").Value = 348100
$ws.Range("J48").Value = 665100
$ws.Range("D49").Value = 1043000
$ws.Range("E49").Value = 935900
$ws.Range("F49").Value = 868200
$ws.Range("G49").Value = 861400
$ws.Range("H49").Value = 816500
$ws.Range("I49").Value = 706200
$ws.Range("J49").Value = 842900
$ws.Range("D52").Value = 88900
$ws.Range("E52").Value = 186700
$ws.Range("F52").Value = 194600
$ws.Range("G52").Value = 174300
$ws.Range("H52").Value = 169100
$ws.Range("I52").Value = 189300
$ws.Range("D54").Value = 2016100
$ws.Range("E54").Value = 2000400
$ws.Range("F54").Value = 1893700
$ws.Range("G54").Value = 1697700
$ws.Range("H54").Value = 1585700
$ws.Range("I54").Value = 1458600
$ws.Range("J54").Value = 1235400
$ws.Range("D57").Value = 349500
$ws.Range("E57").Value = 248200
$ws.Range("F57").Value = 260400
$ws.Range("G57").Value = 171200
$ws.Range("H57").Value = 133300
$ws.Range("I57").Value = 102300
$ws.Range("J57").Value = 109100
$ws.Range("D58").Value = 11800
$ws.Range("H58").Value = 19600
$ws.Range("I58").Value = 15300
$ws.Range("D59").Value = 240500
$ws.Range("E59").Value = 279700
$ws.Range("F59").Value = 258900
$ws.Range("G59").Value = 250000
$ws.Range("H59").Value = 240600
$ws.Range("I59").Value = 192600
$ws.Range("J59").Value = 189300
$ws.Range("D60").Value = 601900
$ws.Range("E60").Value = 535400
$ws.Range("F60").Value = 526600
$ws.Range("G60").Value = 421900
$ws.Range("H60").Value = 393500
$ws.Range("I60").Value = 310300
$ws.Range("J60").Value = 299200
$ws.Range("D61").Value = 634800
$ws.Range("E61").Value = 649100
$ws.Range("F61").Value = 632000
$ws.Range("G61").Value = 535400
$ws.Range("H61").Value = 426300
$ws.Range("I61").Value = 492300
$ws.Range("J61").Value = 549400
$ws.Range("D62").Value = 222400
$ws.Range("E62").Value = 189700
$ws.Range("F62").Value = 137400
$ws.Range("G62").Value = 117100
$ws.Range("H62").Value = 115500
$ws.Range("I62").Value = 71700
$ws.Range("J62").Value = 80100
$ws.Range("D66").Value = 1459100
$ws.Range("E66").Value = 1374200
$ws.Range("F66").Value = 1296000
$ws.Range("G66").Value = 1074300
$ws.Range("H66").Value = 935300
$ws.Range("I66").Value = 874300
$ws.Range("J66").Value = 926100
$ws.Range("D72").Value = 503200
$ws.Range("E72").Value = 577100
$ws.Range("F72").Value = 550900
$ws.Range("G72").Value = 527400
$ws.Range("H72").Value = 526800
$ws.Range("I72").Value = 494100
$ws.Range("J72").Value = 308300
$ws.Range("D76").Value = 557100
$ws.Range("E76").Value = 626200
$ws.Range("F76").Value = 597700
$ws.Range("G76").Value = 623400
$ws.Range("H76").Value = 650400
$ws.Range("I76").Value = 584300
$ws.Range("J76").Value = 309300
$ws.Range("D81").Value = -48600
$ws.Range("E81").Value = 51100
$ws.Range("F81").Value = 47400
$ws.Range("H81").Value = 31100
$ws.Range("I81").Value = 190800
$ws.Range("J81").Value = 30000
$ws.Range("D83").Value = 194500
$ws.Range("E83").Value = 175600
$ws.Range("F83").Value = 127900
$ws.Range("G83").Value = 143500
$ws.Range("H83").Value = 83400
$ws.Range("I83").Value = 64500
$ws.Range("J83").Value = 59700
$ws.Range("D89").Value = 124200
$ws.Range("E89").Value = 296900
$ws.Range("F89").Value = 242700
$ws.Range("G89").Value = 187700
$ws.Range("H89").Value = 94900
$ws.Range("I89").Value = 74600
$ws.Range("J89").Value = 40300
$ws.Range("D91").Value = -17300
$ws.Range("E91").Value = -10900
$ws.Range("F91").Value = -6700
$ws.Range("G91").Value = -7100
$ws.Range("H91").Value = -14100
$ws.Range("I91").Value = -16300
$ws.Range("J91").Value = -20000
$ws.Range("D94").Value = -158200
$ws.Range("E94").Value = -196900
$ws.Range("F94").Value = -136700
$ws.Range("G94").Value = -133400
$ws.Range("H94").Value = -117000
$ws.Range("I94").Value = -63700
$ws.Range("J94").Value = -94200
$ws.Range("D96").Value = -28700
$ws.Range("E96").Value = -30400
$ws.Range("F96").Value = -26200
$ws.Range("J96").Value = -13000
$ws.Range("D100").Value = -29200
$ws.Range("E100").Value = -30900
$ws.Range("F100").Value = -26700
$ws.Range("G100").Value = 58200
$ws.Range("H100").Value = -6500
$ws.Range("I100").Value = 21000
$ws.Range("J100").Value = -50600
$ws.Range("E101").Value = 10500
$ws.Range("F101").Value = 16500
$ws.Range("H101").Value = -8000
$ws.Range("D102").Value = -62900
$ws.Range("E102").Value = 79700
$ws.Range("F102").Value = 95800
$ws.Range("G102").Value = 116600
$ws.Range("H102").Value = -36600
$ws.Range("I102").Value = 31100
$ws.Range("J102").Value = -104400

